# Auto-generated edit script for İş Takip workbook update
# Commit: İş Takip Güncellemesi - 10.01.2026 23:54:27
$wb = $excel.ActiveWorkbook

function Set-TextCell {
    param($ws, $addr, $val)
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws1 = $wb.Worksheets.Item("İş Takip Listesi")

Set-TextCell $ws1 'J2' '2025-06-01'
Set-TextCell $ws1 'K2' '2025-11-01'
Set-TextCell $ws1 'J3' '2025-06-01'
Set-TextCell $ws1 'K3' '2025-11-01'
Set-TextCell $ws1 'J4' '2025-06-01'
Set-TextCell $ws1 'K4' '2025-11-01'
Set-TextCell $ws1 'J5' '2025-06-01'
Set-TextCell $ws1 'K5' '2025-11-01'
Set-TextCell $ws1 'J6' '2025-06-01'
Set-TextCell $ws1 'K6' '2025-11-01'
Set-TextCell $ws1 'J7' '2025-06-01'
Set-TextCell $ws1 'K7' '2025-11-01'
Set-TextCell $ws1 'J8' '2025-06-01'
Set-TextCell $ws1 'K8' '2025-11-01'
Set-TextCell $ws1 'J9' '2025-06-01'
Set-TextCell $ws1 'K9' '2025-11-01'
Set-TextCell $ws1 'J10' '2025-06-01'
Set-TextCell $ws1 'K10' '2025-11-01'
Set-TextCell $ws1 'J33' '2025-06-03'
Set-TextCell $ws1 'K33' '2025-11-03'
Set-TextCell $ws1 'J34' '2025-06-03'
Set-TextCell $ws1 'K34' '2025-11-03'
Set-TextCell $ws1 'J35' '2025-06-03'
Set-TextCell $ws1 'K35' '2025-11-03'
Set-TextCell $ws1 'J36' '2025-06-03'
Set-TextCell $ws1 'K36' '2025-11-03'
Set-TextCell $ws1 'J37' '2025-06-03'
Set-TextCell $ws1 'K37' '2025-11-03'
Set-TextCell $ws1 'J38' '2025-06-03'
Set-TextCell $ws1 'K38' '2025-11-03'
Set-TextCell $ws1 'J39' '2025-06-03'
Set-TextCell $ws1 'K39' '2025-11-03'
Set-TextCell $ws1 'J40' '2025-06-03'
Set-TextCell $ws1 'K40' '2025-11-03'
Set-TextCell $ws1 'J41' '2025-06-03'
Set-TextCell $ws1 'K41' '2025-11-03'
Set-TextCell $ws1 'J42' '2025-06-03'
Set-TextCell $ws1 'K42' '2025-11-03'
Set-TextCell $ws1 'J43' '2025-06-03'
Set-TextCell $ws1 'K43' '2025-11-03'
Set-TextCell $ws1 'J44' '2025-06-03'
Set-TextCell $ws1 'K44' '2025-11-03'
Set-TextCell $ws1 'J45' '2025-06-03'
Set-TextCell $ws1 'K45' '2025-11-03'
Set-TextCell $ws1 'J46' '2025-06-03'
Set-TextCell $ws1 'K46' '2025-11-03'
Set-TextCell $ws1 'J47' '2025-06-03'
Set-TextCell $ws1 'K47' '2025-11-03'
Set-TextCell $ws1 'J48' '2025-06-03'
Set-TextCell $ws1 'K48' '2025-11-03'
Set-TextCell $ws1 'J49' '2025-06-03'
Set-TextCell $ws1 'K49' '2025-11-03'
Set-TextCell $ws1 'J50' '2025-06-03'
Set-TextCell $ws1 'K50' '2025-11-03'
Set-TextCell $ws1 'J51' '2025-06-03'
Set-TextCell $ws1 'K51' '2025-11-03'
Set-TextCell $ws1 'J52' '2025-06-03'
Set-TextCell $ws1 'K52' '2025-11-03'
Set-TextCell $ws1 'J53' '2025-06-03'
Set-TextCell $ws1 'K53' '2025-11-03'
Set-TextCell $ws1 'J54' '2025-06-03'
Set-TextCell $ws1 'K54' '2025-11-03'
Set-TextCell $ws1 'J55' '2025-06-03'
Set-TextCell $ws1 'K55' '2025-11-03'
Set-TextCell $ws1 'J56' '2025-06-03'
Set-TextCell $ws1 'K56' '2025-11-03'
Set-TextCell $ws1 'J57' '2025-06-03'
Set-TextCell $ws1 'K57' '2025-11-03'
Set-TextCell $ws1 'J58' '2025-06-03'
Set-TextCell $ws1 'K58' '2025-11-03'
Set-TextCell $ws1 'J59' '2025-06-03'
Set-TextCell $ws1 'K59' '2025-11-03'
Set-TextCell $ws1 'J60' '2025-06-03'
Set-TextCell $ws1 'K60' '2025-11-03'
Set-TextCell $ws1 'J61' '2025-06-03'
Set-TextCell $ws1 'K61' '2025-11-03'
Set-TextCell $ws1 'J62' '2025-06-03'
Set-TextCell $ws1 'K62' '2025-11-03'
Set-TextCell $ws1 'J63' '2025-06-03'
Set-TextCell $ws1 'K63' '2025-11-03'
Set-TextCell $ws1 'J64' '2025-06-03'
Set-TextCell $ws1 'K64' '2025-11-03'
Set-TextCell $ws1 'J65' '2025-06-03'
Set-TextCell $ws1 'K65' '2025-11-03'
Set-TextCell $ws1 'J66' '2025-06-03'
Set-TextCell $ws1 'K66' '2025-11-03'
Set-TextCell $ws1 'J67' '2025-06-03'
Set-TextCell $ws1 'K67' '2025-11-03'
Set-TextCell $ws1 'J68' '2025-06-03'
Set-TextCell $ws1 'K68' '2025-11-03'
Set-TextCell $ws1 'J69' '2025-06-03'
Set-TextCell $ws1 'K69' '2025-11-03'
Set-TextCell $ws1 'J70' '2025-06-03'
Set-TextCell $ws1 'K70' '2025-11-03'
Set-TextCell $ws1 'J71' '2025-06-03'
Set-TextCell $ws1 'K71' '2025-11-03'
Set-TextCell $ws1 'J72' '2025-06-03'
Set-TextCell $ws1 'K72' '2025-11-03'
Set-TextCell $ws1 'J73' '2025-06-03'
Set-TextCell $ws1 'K73' '2025-11-03'
Set-TextCell $ws1 'J74' '2025-06-03'
Set-TextCell $ws1 'K74' '2025-11-03'
Set-TextCell $ws1 'J75' '2025-06-03'
Set-TextCell $ws1 'K75' '2025-11-03'
Set-TextCell $ws1 'J76' '2025-06-03'
Set-TextCell $ws1 'K76' '2025-11-03'
Set-TextCell $ws1 'J77' '2025-06-03'
Set-TextCell $ws1 'K77' '2025-11-03'
Set-TextCell $ws1 'J78' '2025-06-03'
Set-TextCell $ws1 'K78' '2025-11-03'
Set-TextCell $ws1 'J79' '2025-06-03'
Set-TextCell $ws1 'K79' '2025-11-03'
Set-TextCell $ws1 'J80' '2025-06-03'
Set-TextCell $ws1 'K80' '2025-11-03'
Set-TextCell $ws1 'J81' '2025-06-03'
Set-TextCell $ws1 'K81' '2025-11-03'
Set-TextCell $ws1 'J82' '2025-06-03'
Set-TextCell $ws1 'K82' '2025-11-03'
Set-TextCell $ws1 'J83' '2025-06-03'
Set-TextCell $ws1 'K83' '2025-11-03'
Set-TextCell $ws1 'J84' '2025-06-03'
Set-TextCell $ws1 'K84' '2025-11-03'
Set-TextCell $ws1 'J85' '2025-06-03'
Set-TextCell $ws1 'K85' '2025-11-03'
Set-TextCell $ws1 'J86' '2025-06-03'
Set-TextCell $ws1 'K86' '2025-11-03'
Set-TextCell $ws1 'J87' '2025-06-03'
Set-TextCell $ws1 'K87' '2025-11-03'
Set-TextCell $ws1 'J88' '2025-06-03'
Set-TextCell $ws1 'K88' '2025-11-03'
Set-TextCell $ws1 'J89' '2025-06-03'
Set-TextCell $ws1 'K89' '2025-11-03'
Set-TextCell $ws1 'J90' '2025-06-03'
Set-TextCell $ws1 'K90' '2025-11-03'
Set-TextCell $ws1 'J91' '2025-06-03'
Set-TextCell $ws1 'K91' '2025-11-03'
Set-TextCell $ws1 'J92' '2025-06-03'
Set-TextCell $ws1 'K92' '2025-11-03'
Set-TextCell $ws1 'J93' '2025-06-03'
Set-TextCell $ws1 'K93' '2025-11-03'
Set-TextCell $ws1 'J94' '2025-06-03'
Set-TextCell $ws1 'K94' '2025-11-03'
Set-TextCell $ws1 'J95' '2024-04-01'
Set-TextCell $ws1 'K95' '2025-05-26'
Set-TextCell $ws1 'J96' '2024-04-01'
Set-TextCell $ws1 'K96' '2025-05-26'
Set-TextCell $ws1 'J97' '2024-04-01'
Set-TextCell $ws1 'K97' '2025-05-26'
Set-TextCell $ws1 'J98' '2024-04-01'
Set-TextCell $ws1 'K98' '2025-05-26'
Set-TextCell $ws1 'J99' '2024-04-01'
Set-TextCell $ws1 'K99' '2025-05-26'
Set-TextCell $ws1 'J100' '2024-04-01'
Set-TextCell $ws1 'K100' '2025-05-26'
Set-TextCell $ws1 'J101' '2024-04-01'
Set-TextCell $ws1 'K101' '2025-05-26'
Set-TextCell $ws1 'J102' '2024-04-01'
Set-TextCell $ws1 'K102' '2025-05-26'
Set-TextCell $ws1 'J103' '2024-04-01'
Set-TextCell $ws1 'K103' '2025-05-26'
Set-TextCell $ws1 'J104' '2024-04-01'
Set-TextCell $ws1 'K104' '2025-05-26'
Set-TextCell $ws1 'J105' '2024-04-01'
Set-TextCell $ws1 'K105' '2025-05-26'
Set-TextCell $ws1 'J106' '2024-04-01'
Set-TextCell $ws1 'K106' '2025-05-26'
Set-TextCell $ws1 'J107' '2024-04-01'
Set-TextCell $ws1 'K107' '2025-05-26'
Set-TextCell $ws1 'J108' '2024-04-01'
Set-TextCell $ws1 'K108' '2025-05-26'
Set-TextCell $ws1 'J109' '2024-04-01'
Set-TextCell $ws1 'K109' '2025-05-26'
Set-TextCell $ws1 'J110' '2024-04-01'
Set-TextCell $ws1 'K110' '2025-05-26'
Set-TextCell $ws1 'J111' '2024-04-01'
Set-TextCell $ws1 'K111' '2025-05-26'
Set-TextCell $ws1 'J112' '2024-04-01'
Set-TextCell $ws1 'K112' '2025-05-26'
Set-TextCell $ws1 'J113' '2024-04-01'
Set-TextCell $ws1 'K113' '2025-05-26'
Set-TextCell $ws1 'J114' '2024-04-01'
Set-TextCell $ws1 'K114' '2025-05-26'
Set-TextCell $ws1 'J115' '2024-04-01'
Set-TextCell $ws1 'K115' '2025-05-26'
Set-TextCell $ws1 'J116' '2024-04-01'
Set-TextCell $ws1 'K116' '2025-05-26'
Set-TextCell $ws1 'J117' '2024-04-01'
Set-TextCell $ws1 'K117' '2025-05-26'
Set-TextCell $ws1 'J118' '2024-04-01'
Set-TextCell $ws1 'K118' '2025-05-26'
Set-TextCell $ws1 'J119' '2024-04-01'
Set-TextCell $ws1 'K119' '2025-05-26'
Set-TextCell $ws1 'J120' '2024-04-01'
Set-TextCell $ws1 'K120' '2025-05-26'
Set-TextCell $ws1 'J121' '2024-04-01'
Set-TextCell $ws1 'K121' '2025-05-26'
Set-TextCell $ws1 'J122' '2024-04-01'
Set-TextCell $ws1 'K122' '2025-05-26'
Set-TextCell $ws1 'L114' 'FİRMAYA TESLİM EDİLDİ'

$ws2 = $wb.Worksheets.Item("Güncelleme")

Set-TextCell $ws2 'J2' '2024-07-08'
Set-TextCell $ws2 'N2' '2025-03-11'
Set-TextCell $ws2 'P2' '2025-05-28'
Set-TextCell $ws2 'J3' '2024-10-09'
Set-TextCell $ws2 'N3' '2025-06-29'
Set-TextCell $ws2 'P3' '2025-10-07'
Set-TextCell $ws2 'J4' '2024-08-13'
Set-TextCell $ws2 'N4' '2025-02-03'
Set-TextCell $ws2 'P4' '2025-04-30'
Set-TextCell $ws2 'I5' '2025-02-06'
Set-TextCell $ws2 'J6' '2025-09-19'
Set-TextCell $ws2 'N6' '2025-06-08'
Set-TextCell $ws2 'P6' '2025-12-04'
Set-TextCell $ws2 'I7' '2024-10-09'
Set-TextCell $ws2 'J7' '2024-10-09'
Set-TextCell $ws2 'J8' '2024-09-27'
Set-TextCell $ws2 'N8' '2025-02-26'
Set-TextCell $ws2 'P8' '2025-03-31'
Set-TextCell $ws2 'I9' '2025-05-25'
Set-TextCell $ws2 'J9' '2024-11-10'
Set-TextCell $ws2 'J10' '2024-09-08'
Set-TextCell $ws2 'N10' '2025-06-19'
Set-TextCell $ws2 'P10' '2025-10-11'
Set-TextCell $ws2 'I11' '2025-03-16'
Set-TextCell $ws2 'J11' '2024-10-22'
Set-TextCell $ws2 'N11' '2025-07-09'
Set-TextCell $ws2 'P11' '2025-12-04'
Set-TextCell $ws2 'J12' '2024-09-19'
Set-TextCell $ws2 'N12' '2025-05-29'
Set-TextCell $ws2 'P12' '2025-10-01'
Set-TextCell $ws2 'J13' '2024-11-17'
Set-TextCell $ws2 'J14' '2025-09-15'
Set-TextCell $ws2 'N14' '2025-10-07'
Set-TextCell $ws2 'J15' '2024-12-06'
Set-TextCell $ws2 'N15' '2025-06-26'
Set-TextCell $ws2 'P15' '2025-10-08'
Set-TextCell $ws2 'J16' '2024-08-04'
Set-TextCell $ws2 'N16' '2025-01-12'
Set-TextCell $ws2 'P16' '2025-03-31'
Set-TextCell $ws2 'J17' '2024-08-20'
Set-TextCell $ws2 'N17' '2025-10-07'
Set-TextCell $ws2 'J18' '2025-01-27'
Set-TextCell $ws2 'I19' '2025-03-17'
Set-TextCell $ws2 'J19' '2024-12-06'
Set-TextCell $ws2 'N19' '2025-07-16'
Set-TextCell $ws2 'J20' '2024-11-17'
Set-TextCell $ws2 'N20' '2025-11-21'
Set-TextCell $ws2 'J21' '2024-09-10'
Set-TextCell $ws2 'J22' '2024-09-10'
Set-TextCell $ws2 'J23' '2024-11-18'
Set-TextCell $ws2 'I24' '2025-05-15'
Set-TextCell $ws2 'J25' '2024-10-13'
Set-TextCell $ws2 'J27' '2025-01-03'
Set-TextCell $ws2 'J28' '2024-11-01'
Set-TextCell $ws2 'N28' '2025-10-17'
Set-TextCell $ws2 'I29' '2025-01-21'
Set-TextCell $ws2 'J29' '2024-11-18'
Set-TextCell $ws2 'N29' '2025-09-28'
